# Apply translations-related changes:
#  - survey!F1:   "display.text"  -> "display.prompt.text"
#  - settings!C1: "display.title" -> "display.title.text"
# Plus restore the cursor/selection state that was captured at save time:
#  - settings!  selection moves to C2
#  - survey!    selection moves to F2, and "survey" becomes the active tab
#    (previously "properties" was active)

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Update the translation key text values.
$survey.Range("F1").Value   = "display.prompt.text"
$settings.Range("C1").Value = "display.title.text"

# Restore per-sheet selections.
$settings.Activate() | Out-Null
$settings.Range("C2").Select() | Out-Null

$survey.Activate() | Out-Null
$survey.Range("F2").Select() | Out-Null

Write-Host "done"
